$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 87.667552
$ws.Cells.Item(2, 8).Value = 263.002656
$ws.Cells.Item(2, 9).Value = 0.3606416352150456
$ws.Cells.Item(2, 10).Value = 0.3606416352150456
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 18.04537966666667
$ws.Cells.Item(2, 14).Value = 54.13613900000001
$ws.Cells.Item(2, 15).Value = 0.6797959733292525
$ws.Cells.Item(2, 16).Value = 0.6797959733292525
$ws.Cells.Item(2, 17).Value = 1581.994260287243
$ws.Cells.Item(2, 18).Value = 14237.94834258519
$ws.Cells.Item(2, 19).Value = 0.2451627314340651
$ws.Cells.Item(2, 20).Value = 0.2451627314340651

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 87.667552
$ws.Cells.Item(3, 8).Value = 263.002656
$ws.Cells.Item(3, 9).Value = 0.3606416352150456
$ws.Cells.Item(3, 10).Value = 0.3606416352150456
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.6001993333333334
$ws.Cells.Item(3, 14).Value = 1.800598
$ws.Cells.Item(3, 15).Value = 0.02261039099934159
$ws.Cells.Item(3, 16).Value = 0.02261039099934159
$ws.Cells.Item(3, 17).Value = 52.61800626536534
$ws.Cells.Item(3, 18).Value = 473.5620563882881
$ws.Cells.Item(3, 19).Value = 0.008154248382854102
$ws.Cells.Item(3, 20).Value = 0.008154248382854102

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 87.667552
$ws.Cells.Item(4, 8).Value = 263.002656
$ws.Cells.Item(4, 9).Value = 0.3606416352150456
$ws.Cells.Item(4, 10).Value = 0.3606416352150456
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.690054333333333
$ws.Cells.Item(4, 14).Value = 14.070163
$ws.Cells.Item(4, 15).Value = 0.1766812397072912
$ws.Cells.Item(4, 16).Value = 0.1766812397072912
$ws.Cells.Item(4, 17).Value = 411.1655821503253
$ws.Cells.Item(4, 18).Value = 3700.490239352928
$ws.Cells.Item(4, 19).Value = 0.06371861119985893
$ws.Cells.Item(4, 20).Value = 0.06371861119985894

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 87.667552
$ws.Cells.Item(5, 8).Value = 263.002656
$ws.Cells.Item(5, 9).Value = 0.3606416352150456
$ws.Cells.Item(5, 10).Value = 0.3606416352150456
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.209654333333333
$ws.Cells.Item(5, 14).Value = 9.628962999999999
$ws.Cells.Item(5, 15).Value = 0.1209123959641148
$ws.Cells.Item(5, 16).Value = 0.1209123959641148
$ws.Cells.Item(5, 17).Value = 281.3825381695253
$ws.Cells.Item(5, 18).Value = 2532.442843525728
$ws.Cells.Item(5, 19).Value = 0.04360604419826744
$ws.Cells.Item(5, 20).Value = 0.04360604419826744

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 72.97955566666666
$ws.Cells.Item(6, 8).Value = 218.938667
$ws.Cells.Item(6, 9).Value = 0.3002190170987564
$ws.Cells.Item(6, 10).Value = 0.3002190170987564
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 18.04537966666667
$ws.Cells.Item(6, 14).Value = 54.13613900000001
$ws.Cells.Item(6, 15).Value = 0.6797959733292525
$ws.Cells.Item(6, 16).Value = 0.6797959733292525
$ws.Cells.Item(6, 17).Value = 1316.943789909635
$ws.Cells.Item(6, 18).Value = 11852.49410918671
$ws.Cells.Item(6, 19).Value = 0.2040876789406006
$ws.Cells.Item(6, 20).Value = 0.2040876789406006

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 72.97955566666666
$ws.Cells.Item(7, 8).Value = 218.938667
$ws.Cells.Item(7, 9).Value = 0.3002190170987564
$ws.Cells.Item(7, 10).Value = 0.3002190170987564
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.6001993333333334
$ws.Cells.Item(7, 14).Value = 1.800598
$ws.Cells.Item(7, 15).Value = 0.02261039099934159
$ws.Cells.Item(7, 16).Value = 0.02261039099934159
$ws.Cells.Item(7, 17).Value = 43.80228065809622
$ws.Cells.Item(7, 18).Value = 394.220525922866
$ws.Cells.Item(7, 19).Value = 0.006788069362040901
$ws.Cells.Item(7, 20).Value = 0.006788069362040901

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 72.97955566666666
$ws.Cells.Item(8, 8).Value = 218.938667
$ws.Cells.Item(8, 9).Value = 0.3002190170987564
$ws.Cells.Item(8, 10).Value = 0.3002190170987564
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 4.690054333333333
$ws.Cells.Item(8, 14).Value = 14.070163
$ws.Cells.Item(8, 15).Value = 0.1766812397072912
$ws.Cells.Item(8, 16).Value = 0.1766812397072912
$ws.Cells.Item(8, 17).Value = 342.2780812991912
$ws.Cells.Item(8, 18).Value = 3080.502731692721
$ws.Cells.Item(8, 19).Value = 0.05304306812471272
$ws.Cells.Item(8, 20).Value = 0.05304306812471273

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 72.97955566666666
$ws.Cells.Item(9, 8).Value = 218.938667
$ws.Cells.Item(9, 9).Value = 0.3002190170987564
$ws.Cells.Item(9, 10).Value = 0.3002190170987564
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 3.209654333333333
$ws.Cells.Item(9, 14).Value = 9.628962999999999
$ws.Cells.Item(9, 15).Value = 0.1209123959641148
$ws.Cells.Item(9, 16).Value = 0.1209123959641148
$ws.Cells.Item(9, 17).Value = 234.2391470902578
$ws.Cells.Item(9, 18).Value = 2108.15232381232
$ws.Cells.Item(9, 19).Value = 0.03630020067140218
$ws.Cells.Item(9, 20).Value = 0.03630020067140218

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 41.37117366666666
$ws.Cells.Item(10, 8).Value = 124.113521
$ws.Cells.Item(10, 9).Value = 0.1701903085181653
$ws.Cells.Item(10, 10).Value = 0.1701903085181653
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 18.04537966666667
$ws.Cells.Item(10, 14).Value = 54.13613900000001
$ws.Cells.Item(10, 15).Value = 0.6797959733292525
$ws.Cells.Item(10, 16).Value = 0.6797959733292525
$ws.Cells.Item(10, 17).Value = 746.5585360706021
$ws.Cells.Item(10, 18).Value = 6719.02682463542
$ws.Cells.Item(10, 19).Value = 0.1156946864303119
$ws.Cells.Item(10, 20).Value = 0.1156946864303119

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 41.37117366666666
$ws.Cells.Item(11, 8).Value = 124.113521
$ws.Cells.Item(11, 9).Value = 0.1701903085181653
$ws.Cells.Item(11, 10).Value = 0.1701903085181653
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.6001993333333334
$ws.Cells.Item(11, 14).Value = 1.800598
$ws.Cells.Item(11, 15).Value = 0.02261039099934159
$ws.Cells.Item(11, 16).Value = 0.02261039099934159
$ws.Cells.Item(11, 17).Value = 24.83095085395089
$ws.Cells.Item(11, 18).Value = 223.478557685558
$ws.Cells.Item(11, 19).Value = 0.003848069419894294
$ws.Cells.Item(11, 20).Value = 0.003848069419894294

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 41.37117366666666
$ws.Cells.Item(12, 8).Value = 124.113521
$ws.Cells.Item(12, 9).Value = 0.1701903085181653
$ws.Cells.Item(12, 10).Value = 0.1701903085181653
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 4.690054333333333
$ws.Cells.Item(12, 14).Value = 14.070163
$ws.Cells.Item(12, 15).Value = 0.1766812397072912
$ws.Cells.Item(12, 16).Value = 0.1766812397072912
$ws.Cells.Item(12, 17).Value = 194.0330523304359
$ws.Cells.Item(12, 18).Value = 1746.297470973923
$ws.Cells.Item(12, 19).Value = 0.0300694346951558
$ws.Cells.Item(12, 20).Value = 0.0300694346951558

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 41.37117366666666
$ws.Cells.Item(13, 8).Value = 124.113521
$ws.Cells.Item(13, 9).Value = 0.1701903085181653
$ws.Cells.Item(13, 10).Value = 0.1701903085181653
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 3.209654333333333
$ws.Cells.Item(13, 14).Value = 9.628962999999999
$ws.Cells.Item(13, 15).Value = 0.1209123959641148
$ws.Cells.Item(13, 16).Value = 0.1209123959641148
$ws.Cells.Item(13, 17).Value = 132.7871668343025
$ws.Cells.Item(13, 18).Value = 1195.084501508723
$ws.Cells.Item(13, 19).Value = 0.02057811797280326
$ws.Cells.Item(13, 20).Value = 0.02057811797280326

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 41.06943633333334
$ws.Cells.Item(14, 8).Value = 123.208309
$ws.Cells.Item(14, 9).Value = 0.1689490391680327
$ws.Cells.Item(14, 10).Value = 0.1689490391680327
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 18.04537966666667
$ws.Cells.Item(14, 14).Value = 54.13613900000001
$ws.Cells.Item(14, 15).Value = 0.6797959733292525
$ws.Cells.Item(14, 16).Value = 0.6797959733292525
$ws.Cells.Item(14, 17).Value = 741.1135713309947
$ws.Cells.Item(14, 18).Value = 6670.022141978952
$ws.Cells.Item(14, 19).Value = 0.1148508765242748
$ws.Cells.Item(14, 20).Value = 0.1148508765242748

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 41.06943633333334
$ws.Cells.Item(15, 8).Value = 123.208309
$ws.Cells.Item(15, 9).Value = 0.1689490391680327
$ws.Cells.Item(15, 10).Value = 0.1689490391680327
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.6001993333333334
$ws.Cells.Item(15, 14).Value = 1.800598
$ws.Cells.Item(15, 15).Value = 0.02261039099934159
$ws.Cells.Item(15, 16).Value = 0.02261039099934159
$ws.Cells.Item(15, 17).Value = 24.64984830764245
$ws.Cells.Item(15, 18).Value = 221.848634768782
$ws.Cells.Item(15, 19).Value = 0.003820003834552296
$ws.Cells.Item(15, 20).Value = 0.003820003834552297

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 41.06943633333334
$ws.Cells.Item(16, 8).Value = 123.208309
$ws.Cells.Item(16, 9).Value = 0.1689490391680327
$ws.Cells.Item(16, 10).Value = 0.1689490391680327
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 4.690054333333333
$ws.Cells.Item(16, 14).Value = 14.070163
$ws.Cells.Item(16, 15).Value = 0.1766812397072912
$ws.Cells.Item(16, 16).Value = 0.1766812397072912
$ws.Cells.Item(16, 17).Value = 192.6178878427075
$ws.Cells.Item(16, 18).Value = 1733.560990584367
$ws.Cells.Item(16, 19).Value = 0.0298501256875637
$ws.Cells.Item(16, 20).Value = 0.02985012568756371

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 41.06943633333334
$ws.Cells.Item(17, 8).Value = 123.208309
$ws.Cells.Item(17, 9).Value = 0.1689490391680327
$ws.Cells.Item(17, 10).Value = 0.1689490391680327
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 3.209654333333333
$ws.Cells.Item(17, 14).Value = 9.628962999999999
$ws.Cells.Item(17, 15).Value = 0.1209123959641148
$ws.Cells.Item(17, 16).Value = 0.1209123959641148
$ws.Cells.Item(17, 17).Value = 131.8186942948408
$ws.Cells.Item(17, 18).Value = 1186.368248653567
$ws.Cells.Item(17, 19).Value = 0.02042803312164191
$ws.Cells.Item(17, 20).Value = 0.02042803312164191
